$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the stale _xlchart.v2.* hidden defined names (chart now tracks data differently)
while ($wb.Names.Count -gt 0) {
    $wb.Names.Item(1).Delete()
}

# Row 1 title
$ws.Range("A1").Value = "Fall Quarter"

# Row 2 headers
$ws.Range("B2").Value = "Time Worked This Week"
$ws.Range("C2").Value = "Total Time Worked"
$ws.Range("D2").Value = "Ideal Time Worked"

# Row 3
$ws.Range("B3").Value = 3
$ws.Range("C3").Formula = "=IF(B3>0,B3,#N/A)"
$ws.Range("D3").Value = 12

# Row 4
$ws.Range("B4").Value = 1
$ws.Range("C4").Formula = "=IF(B4>0,C3+B4,#N/A)"
$ws.Range("D4").Formula = "=D3+12"

# Rows 5-13: C column (shared formula "+") and D column (shared formula "+12")
$ws.Range("C5:C13").Formula = "=IF(B5>0,C4+B5,#N/A)"
$ws.Range("D5:D12").Formula = "=D4+12"

# Row 13 - D13 no longer a formula, just value 120
$ws.Range("D13").Value = 120

# Widen column B to fit the new longer header text
$ws.Columns.Item(2).ColumnWidth = 19.833333333333336

# Update the view: zoom level and active selection
$excel.ActiveWindow.Zoom = 109
[void]$ws.Range("F30").Select()
